$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new entry was logged in the "Dominik Deszczka" section of the table
# (columns H:J, row 9): a date, a file/task description and a line count.
$ws.Range("H9").Value2 = 45770
$ws.Range("I9").Value = "Tworzenie klas"
$ws.Range("J9").Value = 118

# Give the new date cell the same formatting as the other plain date
# cells used throughout the table (e.g. B7/B8: date number format, no
# extra horizontal alignment) by copying their format across.
$ws.Range("B7").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column I now holds a longer text value ("Tworzenie klas"); widen it so
# the text fits, matching what Excel's auto-fit would do after typing it.
$ws.Columns.Item(9).ColumnWidth = 25.6

# Leave the selection where the user ended up after entering the new row.
$ws.Range("I12").Select()
